$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.204118999999999
$ws.Range("H2").Value = 12.612357
$ws.Range("I2").Value = 0.01983154129720676
$ws.Range("J2").Value = 0.01983154129720676
$ws.Range("M2").Value = 40.81054266666667
$ws.Range("N2").Value = 122.431628
$ws.Range("O2").Value = 0.9943414173631485
$ws.Range("P2").Value = 0.9943414173631485
$ws.Range("Q2").Value = 171.572377825244
$ws.Range("R2").Value = 1544.151400427196
$ws.Range("S2").Value = 0.01971932288196038
$ws.Range("T2").Value = 0.01971932288196038

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.204118999999999
$ws.Range("H3").Value = 12.612357
$ws.Range("I3").Value = 0.01983154129720676
$ws.Range("J3").Value = 0.01983154129720676
$ws.Range("O3").Value = 0.0002749163555820933
$ws.Range("P3").Value = 0.0002749163555820933
$ws.Range("Q3").Value = 0.04743647604999999
$ws.Range("R3").Value = 0.42692828445
$ws.Range("S3").Value = 0.00000545201505900386
$ws.Range("T3").Value = 0.000005452015059003862

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.204118999999999
$ws.Range("H4").Value = 12.612357
$ws.Range("I4").Value = 0.01983154129720676
$ws.Range("J4").Value = 0.01983154129720676
$ws.Range("O4").Value = 0.002616675800765965
$ws.Range("P4").Value = 0.002616675800765965
$ws.Range("Q4").Value = 0.451504162751
$ws.Range("R4").Value = 4.063537464758999
$ws.Range("S4").Value = 0.00005189271420429179
$ws.Range("T4").Value = 0.0000518927142042918

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.204118999999999
$ws.Range("H5").Value = 12.612357
$ws.Range("I5").Value = 0.01983154129720676
$ws.Range("J5").Value = 0.01983154129720676
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.113565
$ws.Range("N5").Value = 0.340695
$ws.Range("O5").Value = 0.002766990480503436
$ws.Range("P5").Value = 0.002766990480503436
$ws.Range("Q5").Value = 0.477440774235
$ws.Range("R5").Value = 4.296966968115
$ws.Range("S5").Value = 0.00005487368598308185
$ws.Range("T5").Value = 0.00005487368598308186

# Row 6
$ws.Range("I6").Value = 0.8539093107807857
$ws.Range("J6").Value = 0.8539093107807858
$ws.Range("M6").Value = 40.81054266666667
$ws.Range("N6").Value = 122.431628
$ws.Range("O6").Value = 0.9943414173631485
$ws.Range("P6").Value = 0.9943414173631485
$ws.Range("Q6").Value = 7387.587717068163
$ws.Range("R6").Value = 66488.28945361346
$ws.Range("S6").Value = 0.8490773943813558
$ws.Range("T6").Value = 0.8490773943813559

# Row 7
$ws.Range("I7").Value = 0.8539093107807857
$ws.Range("J7").Value = 0.8539093107807858
$ws.Range("O7").Value = 0.0002749163555820933
$ws.Range("P7").Value = 0.0002749163555820933
$ws.Range("S7").Value = 0.0002347536357174707
$ws.Range("T7").Value = 0.0002347536357174708

# Row 8
$ws.Range("I8").Value = 0.8539093107807857
$ws.Range("J8").Value = 0.8539093107807858
$ws.Range("O8").Value = 0.002616675800765965
$ws.Range("P8").Value = 0.002616675800765965
$ws.Range("S8").Value = 0.002234403829568825
$ws.Range("T8").Value = 0.002234403829568826

# Row 9
$ws.Range("I9").Value = 0.8539093107807857
$ws.Range("J9").Value = 0.8539093107807858
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.113565
$ws.Range("N9").Value = 0.340695
$ws.Range("O9").Value = 0.002766990480503436
$ws.Range("P9").Value = 0.002766990480503436
$ws.Range("Q9").Value = 20.55771240146001
$ws.Range("R9").Value = 185.01941161314
$ws.Range("S9").Value = 0.002362758934143684
$ws.Range("T9").Value = 0.002362758934143684

# Row 10
$ws.Range("G10").Value = 26.057747
$ws.Range("H10").Value = 78.173241
$ws.Range("I10").Value = 0.1229188055196976
$ws.Range("J10").Value = 0.1229188055196976
$ws.Range("M10").Value = 40.81054266666667
$ws.Range("N10").Value = 122.431628
$ws.Range("O10").Value = 0.9943414173631485
$ws.Range("P10").Value = 0.9943414173631485
$ws.Range("Q10").Value = 1063.430795740705
$ws.Range("R10").Value = 9570.877161666349
$ws.Range("S10").Value = 0.1222232593010413
$ws.Range("T10").Value = 0.1222232593010413

# Row 11
$ws.Range("G11").Value = 26.057747
$ws.Range("H11").Value = 78.173241
$ws.Range("I11").Value = 0.1229188055196976
$ws.Range("J11").Value = 0.1229188055196976
$ws.Range("O11").Value = 0.0002749163555820933
$ws.Range("P11").Value = 0.0002749163555820933
$ws.Range("Q11").Value = 0.2940182453166667
$ws.Range("R11").Value = 2.64616420785
$ws.Range("S11").Value = 0.00003379239004597936
$ws.Range("T11").Value = 0.00003379239004597937

# Row 12
$ws.Range("G12").Value = 26.057747
$ws.Range("H12").Value = 78.173241
$ws.Range("I12").Value = 0.1229188055196976
$ws.Range("J12").Value = 0.1229188055196976
$ws.Range("O12").Value = 0.002616675800765965
$ws.Range("P12").Value = 0.002616675800765965
$ws.Range("Q12").Value = 2.798489110896333
$ws.Range("R12").Value = 25.186401998067
$ws.Range("S12").Value = 0.0003216386638624506
$ws.Range("T12").Value = 0.0003216386638624506

# Row 13
$ws.Range("G13").Value = 26.057747
$ws.Range("H13").Value = 78.173241
$ws.Range("I13").Value = 0.1229188055196976
$ws.Range("J13").Value = 0.1229188055196976
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.113565
$ws.Range("N13").Value = 0.340695
$ws.Range("O13").Value = 0.002766990480503436
$ws.Range("P13").Value = 0.002766990480503436
$ws.Range("Q13").Value = 2.959248038055001
$ws.Range("R13").Value = 26.633232342495
$ws.Range("S13").Value = 0.0003401151647478565
$ws.Range("T13").Value = 0.0003401151647478565

# Row 14
$ws.Range("G14").Value = 0.7081243333333332
$ws.Range("H14").Value = 2.124373
$ws.Range("I14").Value = 0.003340342402309973
$ws.Range("J14").Value = 0.003340342402309974
$ws.Range("M14").Value = 40.81054266666667
$ws.Range("N14").Value = 122.431628
$ws.Range("O14").Value = 0.9943414173631485
$ws.Range("P14").Value = 0.9943414173631485
$ws.Range("Q14").Value = 28.89893831880489
$ws.Range("R14").Value = 260.090444869244
$ws.Range("S14").Value = 0.003321440798791123
$ws.Range("T14").Value = 0.003321440798791124

# Row 15
$ws.Range("G15").Value = 0.7081243333333332
$ws.Range("H15").Value = 2.124373
$ws.Range("I15").Value = 0.003340342402309973
$ws.Range("J15").Value = 0.003340342402309974
$ws.Range("O15").Value = 0.0002749163555820933
$ws.Range("P15").Value = 0.0002749163555820933
$ws.Range("Q15").Value = 0.007990002894444443
$ws.Range("R15").Value = 0.07191002604999999
$ws.Range("S15").Value = 0.0000009183147596393923
$ws.Range("T15").Value = 0.0000009183147596393927

# Row 16
$ws.Range("G16").Value = 0.7081243333333332
$ws.Range("H16").Value = 2.124373
$ws.Range("I16").Value = 0.003340342402309973
$ws.Range("J16").Value = 0.003340342402309974
$ws.Range("O16").Value = 0.002616675800765965
$ws.Range("P16").Value = 0.002616675800765965
$ws.Range("Q16").Value = 0.07604948486122222
$ws.Range("R16").Value = 0.684445363751
$ws.Range("S16").Value = 0.000008740593130396954
$ws.Range("T16").Value = 0.000008740593130396958

# Row 17
$ws.Range("G17").Value = 0.7081243333333332
$ws.Range("H17").Value = 2.124373
$ws.Range("I17").Value = 0.003340342402309973
$ws.Range("J17").Value = 0.003340342402309974
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.113565
$ws.Range("N17").Value = 0.340695
$ws.Range("O17").Value = 0.002766990480503436
$ws.Range("P17").Value = 0.002766990480503436
$ws.Range("Q17").Value = 0.08041813991500001
$ws.Range("R17").Value = 0.723763259235
$ws.Range("S17").Value = 0.000009242695628813674
$ws.Range("T17").Value = 0.000009242695628813675
